# Add a "PolynomialCalculator" worksheet after the existing "BasicCalculator"
# sheet, populate it, and refresh a couple of workbook-level properties.

$wb = $excel.ActiveWorkbook

# --- workbook-level bits -------------------------------------------------
# (absolute path recorded by the authoring machine)
$wb.Application.ActiveWorkbook.Path | Out-Null

# --- add the new sheet, right after BasicCalculator ----------------------
$basic = $wb.Worksheets.Item(1)
$poly = $wb.Worksheets.Add($null, $basic)
$poly.Name = "PolynomialCalculator"

# --- header rows (merged I:M, one label each) -----------------------------
$labels = @("input", "table", "simplified", "range")
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 1
    $rng = $poly.Range("I" + $row + ":M" + $row)
    $rng.Merge()
    $rng.Value = $labels[$i]
    $rng.Font.Size = 20
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
}

# --- M7 / M8 placeholder cells (vertical-center only, no text) -----------
foreach ($row in 7, 8) {
    $cell = $poly.Range("M" + $row)
    $cell.Font.Size = 20
    $cell.VerticalAlignment = -4108
}

# --- row heights across the used area (matches BasicCalculator's look) ---
$poly.Range("A1:A28").EntireRow.RowHeight = 48

# --- selection / view state ------------------------------------------------
$poly.Range("I5").Select()

Write-Output "PolynomialCalculator sheet added"
